$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.823.06'
$ws.Range("E2").Value = '  -2.23%  '
$ws.Range("D3").Value = '1.856.35'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = '''304.59'
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("D7").Value = '''0.5072'
$ws.Range("E7").Value = '  -2.75%  '
$ws.Range("D8").Value = '''0.3649'
$ws.Range("E8").Value = '  -3.44%  '
$ws.Range("D9").Value = '''0.07144'
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").Value = '''0.8883'
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '''20.73'
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("D12").Value = '''0.07525'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").Value = '1.858.92'
$ws.Range("E13").Value = '  -1.75%  '
$ws.Range("D14").Value = '''91.33'
$ws.Range("E14").Value = '  +0.65%  '
$ws.Range("D15").Value = '''5.235'
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '''0.000008525'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '26.877.48'
$ws.Range("E20").Value = '  -2.21%  '
$ws.Range("D21").Value = '''5.008'
$ws.Range("E21").Value = '  -2.22%  '
$ws.Range("D22").Value = '2.094.14'
$ws.Range("E22").Value = '  -1.29%  '
$ws.Range("D23").Value = '''10.24'
$ws.Range("E23").Value = '  -4.39%  '
$ws.Range("D24").Value = '''6.431'
$ws.Range("E24").Value = '  -1.73%  '
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").Value = '''146.65'
$ws.Range("E26").Value = '  -3.85%  '
$ws.Range("D27").Value = '''17.81'
$ws.Range("E27").Value = '  -1.89%  '
$ws.Range("D28").Value = '''2.047'
$ws.Range("E28").Value = '  -5.54%  '
$ws.Range("D29").Value = '''113.13'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '''4.637'
$ws.Range("E30").Value = '  -3.16%  '
$ws.Range("D31").Value = '''4.674'
$ws.Range("E31").Value = '  -1.40%  '
$ws.Range("D32").Value = '''0.09243'
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("D33").Value = '''0.05105'
$ws.Range("E33").Value = '  -2.18%  '
$ws.Range("D34").Value = '''3.068'
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").Value = '''1.147'
$ws.Range("E35").Value = '  -4.45%  '
$ws.Range("D36").Value = '''0.7292'
$ws.Range("E36").Value = '  -3.85%  '
$ws.Range("D37").Value = '''3.187'
$ws.Range("E37").Value = '  +3.91%  '
$ws.Range("D38").Value = '''0.02008'
$ws.Range("E38").Value = '  -2.61%  '
$ws.Range("D39").Value = '''2.446'
$ws.Range("E39").Value = '  -3.27%  '
$ws.Range("E40").Value = '  -1.37%  '
$ws.Range("D41").Value = '''0.5275'
$ws.Range("E41").Value = '  -4.12%  '
$ws.Range("D42").Value = '''117.58'
$ws.Range("E42").Value = '  +1.53%  '
$ws.Range("D43").Value = '''6.479'
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").Value = '''8.448'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").Value = '''0.9998'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").Value = '''0.4629'
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("D48").Value = '''9.882'
$ws.Range("E48").Value = '  -4.69%  '
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = '''37.12'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").Value = '''62.85'
$ws.Range("E51").Value = '  -4.50%  '
